$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Timor-Leste (TLS) data: 4 new rows appended after the existing data (rows 155-158).
# Columns: A = date posted, B = base URL for AWS, C = country code,
#          J = Vector Tiles .mbtiles link, K = Dynamic Vectortiles (pbf) link

$startRow = 155
$lastRow = 154
$dateSerial = 44376
$baseUrl = "https://itos-humanitarian.s3.amazonaws.com"
$country = "TLS"

$pbfUrls = @(
    "http://apps.itos.uga.edu/CODV2API/api/v1/Themes/cod-ab/locations/TLS/versions/current/0/{z}/{x}/{y}.pbf",
    "http://apps.itos.uga.edu/CODV2API/api/v1/Themes/cod-ab/locations/TLS/versions/current/1/{z}/{x}/{y}.pbf",
    "http://apps.itos.uga.edu/CODV2API/api/v1/Themes/cod-ab/locations/TLS/versions/current/2/{z}/{x}/{y}.pbf",
    "http://apps.itos.uga.edu/CODV2API/api/v1/Themes/cod-ab/locations/TLS/versions/current/2/{z}/{x}/{y}.pbf"
)

$mbtUrls = @(
    "https://itos-humanitarian.s3.amazonaws.com/v1/VectorTile/COD_TLS/Admin0-MBT/Admin0.mbtiles",
    "https://itos-humanitarian.s3.amazonaws.com/v1/VectorTile/COD_TLS/Admin1-MBT/Admin1.mbtiles",
    "https://itos-humanitarian.s3.amazonaws.com/v1/VectorTile/COD_TLS/Admin2-MBT/Admin2.mbtiles",
    "https://itos-humanitarian.s3.amazonaws.com/v1/VectorTile/COD_TLS/Admin3-MBT/Admin3.mbtiles"
)

# Fill A, B and C columns first (A/C use new values; B reuses an existing shared string)
for ($i = 0; $i -lt 4; $i++) {
    $r = $startRow + $i

    $aCell = $ws.Cells.Item($r, 1)
    $ws.Cells.Item($lastRow, 1).Copy() | Out-Null
    $aCell.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $aCell.Value = $dateSerial

    $ws.Cells.Item($r, 2).Value = $baseUrl
    $ws.Cells.Item($r, 3).Value = $country
}

# Fill column K (pbf / Dynamic Vectortiles links) for all four rows
for ($i = 0; $i -lt 4; $i++) {
    $r = $startRow + $i
    $kCell = $ws.Cells.Item($r, 11)
    $ws.Cells.Item($lastRow, 11).Copy() | Out-Null
    $kCell.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $kCell.Value = $pbfUrls[$i]
    $ws.Hyperlinks.Add($kCell, $pbfUrls[$i]) | Out-Null
    $kCell.PasteSpecial(-4122) | Out-Null   # restore original (non-hyperlink-duplicated) format
}

# Fill column J (Vector Tiles .mbtiles links) for all four rows
for ($i = 0; $i -lt 4; $i++) {
    $r = $startRow + $i
    $jCell = $ws.Cells.Item($r, 10)
    $ws.Cells.Item($lastRow, 10).Copy() | Out-Null
    $jCell.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $jCell.Value = $mbtUrls[$i]
    $ws.Hyperlinks.Add($jCell, $mbtUrls[$i]) | Out-Null
    $jCell.PasteSpecial(-4122) | Out-Null   # restore original (non-hyperlink-duplicated) format
}
